$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 1.53
$ws.Range("G2").Value = 1.6
$ws.Range("H2").Value = 6.2
$ws.Range("O2").Value = 1.27
$ws.Range("P2").Value = 2.12
$ws.Range("Q2").Value = 1.78
$ws.Range("R2").Value = 1.43
$ws.Range("S2").Value = 3.05
$ws.Range("T2").Value = 1.91
$ws.Range("U2").Value = 1.93
$ws.Range("W2").Value = 2.66
$ws.Range("AC2").Value = 11.5
$ws.Range("AD2").Value = 55
$ws.Range("AK2").Value = 32
$ws.Range("F3").Value = 7.2
$ws.Range("G3").Value = 8
$ws.Range("H3").Value = 1.43
$ws.Range("I3").Value = 1.48
$ws.Range("J3").Value = 5.1
$ws.Range("K3").Value = 5.9
$ws.Range("N3").Value = 5.4
$ws.Range("O3").Value = 1.19
$ws.Range("P3").Value = 2.52
$ws.Range("Q3").Value = 1.59
$ws.Range("R3").Value = 1.59
$ws.Range("S3").Value = 2.5
$ws.Range("T3").Value = 1.81
$ws.Range("U3").Value = 2.06
$ws.Range("V3").Value = 3.05
$ws.Range("X3").Value = 30
$ws.Range("AA3").Value = 13.5
$ws.Range("AG3").Value = 27
$ws.Range("AH3").Value = 38
$ws.Range("AI3").Value = 80
$ws.Range("F4").Value = 1.45
$ws.Range("G4").Value = 1.49
$ws.Range("H4").Value = 7
$ws.Range("I4").Value = 19.5
$ws.Range("K4").Value = 5.9
$ws.Range("P4").Value = 2.36
$ws.Range("Q4").Value = 1.66
$ws.Range("V4").Value = 1.1
$ws.Range("W4").Value = 3
$ws.Range("F5").Value = 1.99
$ws.Range("G5").Value = 2.14
$ws.Range("H5").Value = 4
$ws.Range("I5").Value = 4.7
$ws.Range("J5").Value = 3.25
$ws.Range("K5").Value = 3.8
$ws.Range("L5").Value = 1.45
$ws.Range("N5").Value = 3.35
$ws.Range("O5").Value = 1.37
$ws.Range("P5").Value = 1.78
$ws.Range("Q5").Value = 2.12
$ws.Range("R5").Value = 1.29
$ws.Range("S5").Value = 3.85
$ws.Range("T5").Value = 1.87
$ws.Range("U5").Value = 1.93
$ws.Range("V5").Value = 1.28
$ws.Range("W5").Value = 1.87
$ws.Range("Y5").Value = 14.5
$ws.Range("Z5").Value = 46
$ws.Range("AA5").Value = 1000
$ws.Range("AB5").Value = 8.6
$ws.Range("AC5").Value = 8
$ws.Range("AE5").Value = 150
$ws.Range("AF5").Value = 13
$ws.Range("AG5").Value = 11.5
$ws.Range("AH5").Value = 20
$ws.Range("AI5").Value = 1000
$ws.Range("AJ5").Value = 44
$ws.Range("AK5").Value = 25
$ws.Range("AL5").Value = 130
$ws.Range("AN5").Value = 19
$ws.Range("AO5").Value = 1000
$ws.Range("F6").Value = 1.8
$ws.Range("G6").Value = 1.88
$ws.Range("H6").Value = 4.9
$ws.Range("I6").Value = 6.2
$ws.Range("J6").Value = 3.4
$ws.Range("K6").Value = 3.95
$ws.Range("L6").Value = 1.47
$ws.Range("M6").Value = 1.08
$ws.Range("N6").Value = 3.3
$ws.Range("P6").Value = 1.78
$ws.Range("Q6").Value = 2.12
$ws.Range("S6").Value = 4
$ws.Range("T6").Value = 1.95
$ws.Range("U6").Value = 1.87
$ws.Range("X6").Value = 19.5
$ws.Range("Y6").Value = 990
$ws.Range("AB6").Value = 7.8
$ws.Range("AC6").Value = 8.6
$ws.Range("AD6").Value = 60
$ws.Range("AF6").Value = 11
$ws.Range("AG6").Value = 40
$ws.Range("AL6").Value = 290
